$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to
# remain text (matching the source data, which stores prices as strings),
# otherwise Excel auto-converts numeric-looking text into a Number.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D14", "D16", "D18", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D42", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.764.83"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "2.266.38"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "304.26"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").Value = "91.33"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  +3.42%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "32.28"
$ws.Range("E10").Value = "  +4.81%  "
$ws.Range("D11").Value = "53.56"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "6.57"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "2.617.43"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "14.16"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "2.295.41"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").Value = "0.759"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "41.700.02"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("D20").Value = "12.13"
$ws.Range("E20").Value = "  +7.76%  "
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "5.91"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "66.71"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").Value = "241.52"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "1.90"
$ws.Range("E27").Value = "  +4.39%  "
$ws.Range("D28").Value = "24.12"
$ws.Range("E28").Value = "  +3.85%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +10.11%  "
$ws.Range("D30").Value = "9.51"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "159.79"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "34.09"
$ws.Range("E32").Value = "  +7.39%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "5.14"
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").Value = "16.65"
$ws.Range("E38").Value = "  +6.42%  "
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("D42").Value = "3.90"
$ws.Range("E42").Value = "  +3.79%  "
$ws.Range("D43").Value = "2.059.18"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "19.39"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("E48").Value = "  +6.66%  "
$ws.Range("D49").Value = "73.10"
$ws.Range("E49").Value = "  +7.82%  "
$ws.Range("D50").Value = "1.51"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").Value = "1.15"
$ws.Range("E51").Value = "  +2.39%  "
